# Generate Report for handback
# Update the "Correspond Handoff Datetime" (D3) and "Correspond Handback
# DateTime" (G3) timestamps on the 56bec4ee... row of the per-language
# report sheets, for both the zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("D3").Value = "2016-01-22 02:32:50"
$wsZhCn.Range("G3").Value = "2016-01-22 02:33:35"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("D3").Value = "2016-01-22 02:33:05"
$wsDeDe.Range("G3").Value = "2016-01-22 02:33:58"
